$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Row 3: replace TC_ID value and fill in USERNAME/PASSWORD/TRIP_TYPE/SERVICE_CLASS
$ws.Range("A3").Value = "loginWithValidCredentialsViral"
$ws.Range("B3").Value = "mercury"
$ws.Range("C3").Value = "mercury"
$ws.Range("D3").Value = "OneWay"
$ws.Range("E3").Value = "Business"

# New cell C7, center-aligned, empty
$ws.Range("C7").Value = ""
$ws.Range("C7").HorizontalAlignment = -4108  # xlCenter

# Update selection to A4
$ws.Range("A4").Select()
